# Weekly data refresh: a new observation (week) was recorded for
# "Vega Monumental Concepción - Cebollín". It is inserted as the new
# row 97, pushing the previously-existing rows 97:119 down to 98:120.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 97 (shifts rows 97:119 -> 98:120, same as
# Excel's Rows(97).Insert, and carries the column-D date style down).
$ws.Rows("97:97").Insert()

# Populate the newly inserted row 97 with the new weekly record.
$ws.Range("A97").Value = 11
$ws.Range("B97").Value = 'Vega Monumental Concepción'
$ws.Range("C97").Value = 'Bíobío'
$ws.Range("D97").Value = '2023-06-16'
$ws.Range("E97").Value = 8
$ws.Range("F97").Value = 100112037
$ws.Range("G97").Value = 'Cebollín'
$ws.Range("H97").Value = 'Sin especificar'
$ws.Range("I97").Value = 'Primera'
$ws.Range("J97").Value = 100
$ws.Range("K97").Value = 3500
$ws.Range("L97").Value = 4000
$ws.Range("M97").Value = 3750
$ws.Range("N97").Value = '$/paquete 36 unidades'
$ws.Range("O97").Value = 'Región Metropolitana'
$ws.Range("P97").Value = 104
$ws.Range("Q97").Value = 36
$ws.Range("R97").Value = 'Hortaliza'
